# Update workbook "horarios-141-completo" with the 31/12/2025 11:46 scrape.
# Sheets: LP1912 (main feed), LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": 17 new rows (886-902), dimension A1:G885 -> A1:G902,
# header timestamp + row-count refresh.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 11:46:18"
$ws1.Cells.Item(3,1).Value = "Total filas: 901"

$sheet1Rows = @(
    @(886, "11:46:08", "11:50", "11_ETCHEVERRY",   4, "LP1912", "31/12/2025"),
    @(887, "11:46:08", "11:53", "15_ABASTO",        7, "LP1912", "31/12/2025"),
    @(888, "11:46:08", "11:54", "225_GOMEZ",        8, "LP1912", "31/12/2025"),
    @(889, "11:46:08", "11:57", "16_SANTA ANA",    11, "LP1912", "31/12/2025"),
    @(890, "11:46:08", "12:04", "23_HERNANDEZ",    18, "LP1912", "31/12/2025"),
    @(891, "11:46:08", "12:05", "17_ROMERO",       19, "LP1912", "31/12/2025"),
    @(892, "11:46:08", "12:09", "16_SANTA ANA",    23, "LP1912", "31/12/2025"),
    @(893, "11:46:08", "12:17", "15_ABASTO",       31, "LP1912", "31/12/2025"),
    @(894, "11:46:08", "12:18", "10_OLMOS",        32, "LP1912", "31/12/2025"),
    @(895, "11:46:08", "12:18", "17_ROMERO",       32, "LP1912", "31/12/2025"),
    @(896, "11:46:08", "12:29", "215C_EL PATO",    43, "LP1912", "31/12/2025"),
    @(897, "11:46:08", "12:39", "23_HERNANDEZ",    53, "LP1912", "31/12/2025"),
    @(898, "11:46:08", "12:40", "15X38_ABASTO",    54, "LP1912", "31/12/2025"),
    @(899, "11:46:08", "12:51", "15_ABASTO",       65, "LP1912", "31/12/2025"),
    @(900, "11:46:08", "13:04", "23_HERNANDEZ",    78, "LP1912", "31/12/2025"),
    @(901, "11:46:08", "13:07", "14_ABASTO",       81, "LP1912", "31/12/2025"),
    @(902, "11:46:08", "13:21", "17_ROMERO",       95, "LP1912", "31/12/2025")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    # Column A stays blank (as in every other detail row of this sheet).
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": 1 new row (69), dimension A1:G68 -> A1:G69,
# header timestamp + row-count refresh.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 11:46:18"
$ws2.Cells.Item(3,1).Value = "Total filas: 68"

$ws2.Cells.Item(69,2).Value = "31/12/2025"
$ws2.Cells.Item(69,3).Value = "11:46:08"
$ws2.Cells.Item(69,4).Value = "12:29"
$ws2.Cells.Item(69,5).Value = "215C_EL PATO"
$ws2.Cells.Item(69,6).Value = 43
$ws2.Cells.Item(69,7).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": 2 new rows (105-106), dimension A1:G104 -> A1:G106,
# header timestamp + row-count refresh.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 11:46:18"
$ws3.Cells.Item(3,1).Value = "Total filas: 105"

$ws3.Cells.Item(105,2).Value = "31/12/2025"
$ws3.Cells.Item(105,3).Value = "11:46:18"
$ws3.Cells.Item(105,4).Value = "13:09"
$ws3.Cells.Item(105,5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(105,6).Value = 83
$ws3.Cells.Item(105,7).Value = "L6173"

$ws3.Cells.Item(106,2).Value = "31/12/2025"
$ws3.Cells.Item(106,3).Value = "11:46:18"
$ws3.Cells.Item(106,4).Value = "13:14"
$ws3.Cells.Item(106,5).Value = "215A_LA PLATA"
$ws3.Cells.Item(106,6).Value = 88
$ws3.Cells.Item(106,7).Value = "L6173"
